$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row's data (Fecha, Calidad, Volumen, Precio minimo, Precio maximo,
# Precio promedio ponderado, Unidad de comercializacion, Precio $/Kg) is
# being swapped between rows, per the source diff. Row 8 stays unchanged.

$rows = @{
    2  = @{ D = 45084; L = "Primera"; M = 100; N = 20000; O = 21000; P = 20500; Q = "`$/caja 18 kilos granel"; S = 1139 }
    3  = @{ D = 45155; L = "Primera"; M = 40;  N = 25000; O = 26000; P = 25500; Q = "`$/caja 18 kilos";        S = 1417 }
    4  = @{ D = 45168; L = "Primera"; M = 50;  N = 26000; O = 26000; P = 26000; Q = "`$/caja 18 kilos";        S = 1444 }
    5  = @{ D = 45168; L = "Segunda"; M = 50;  N = 22000; O = 22000; P = 22000; Q = "`$/caja 18 kilos";        S = 1222 }
    6  = @{ D = 44699; L = "Primera"; M = 100; N = 20000; O = 22000; P = 21000; Q = "`$/caja 18 kilos";        S = 1167 }
    7  = @{ D = 44699; L = "Segunda"; M = 50;  N = 18000; O = 18000; P = 18000; Q = "`$/caja 18 kilos";        S = 1000 }
    9  = @{ D = 45002; L = "Primera"; M = 100; N = 12000; O = 13000; P = 12500; Q = "`$/caja 18 kilos";        S = 694 }
    10 = @{ D = 45030; L = "Primera"; M = 100; N = 15000; O = 16000; P = 15500; Q = "`$/caja 18 kilos granel"; S = 861 }
    11 = @{ D = 45044; L = "Primera"; M = 100; N = 17000; O = 18000; P = 17500; Q = "`$/caja 18 kilos";        S = 972 }
    12 = @{ D = 44316; L = "Primera"; M = 50;  N = 20000; O = 20000; P = 20000; Q = "`$/caja 18 kilos";        S = 1111 }
    13 = @{ D = 44687; L = "Primera"; M = 100; N = 18000; O = 19000; P = 18500; Q = "`$/caja 18 kilos";        S = 1028 }
    14 = @{ D = 44280; L = "Primera"; M = 100; N = 14000; O = 15000; P = 14500; Q = "`$/caja 18 kilos";        S = 806 }
    15 = @{ D = 44280; L = "Segunda"; M = 50;  N = 12000; O = 12000; P = 12000; Q = "`$/caja 18 kilos";        S = 667 }
    16 = @{ D = 45014; L = "Primera"; M = 50;  N = 13000; O = 14000; P = 13600; Q = "`$/caja 18 kilos";        S = 756 }
    17 = @{ D = 45014; L = "Segunda"; M = 20;  N = 10000; O = 10000; P = 10000; Q = "`$/caja 18 kilos";        S = 556 }
    18 = @{ D = 44819; L = "Primera"; M = 100; N = 25000; O = 26000; P = 25500; Q = "`$/caja 18 kilos granel"; S = 1417 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("Q$r").Value = $vals.Q
    $ws.Range("S$r").Value = $vals.S
}
